$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos list with the latest scraped price/volume figures
# (and the re-ranked rows for Filecoin/Mantle and Bittensor/Maker/Cosmos).

$ws.Cells.Item(2, 4).Value = '70.713.07'
$ws.Cells.Item(2, 5).Value = '  +5.81%  '

$ws.Cells.Item(3, 4).Value = '3.815.04'
$ws.Cells.Item(3, 5).Value = '  +23.80%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '618.64'
$ws.Cells.Item(5, 5).Value = '  +8.48%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '178.77'
$ws.Cells.Item(6, 5).Value = '  +1.87%  '

$ws.Cells.Item(7, 4).Value = '3.813.16'
$ws.Cells.Item(7, 5).Value = '  +23.79%  '

$ws.Cells.Item(8, 5).Value = '  -0.11%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.549'
$ws.Cells.Item(9, 5).Value = '  +7.01%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.173'
$ws.Cells.Item(10, 5).Value = '  +14.36%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.40'
$ws.Cells.Item(11, 5).Value = '  -0.15%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.507'
$ws.Cells.Item(12, 5).Value = '  +8.94%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '40.92'
$ws.Cells.Item(13, 5).Value = '  +14.40%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.0000264'
$ws.Cells.Item(14, 5).Value = '  +10.50%  '

$ws.Cells.Item(15, 4).Value = '4.447.47'
$ws.Cells.Item(15, 5).Value = '  +23.60%  '

$ws.Cells.Item(16, 4).Value = '3.801.88'
$ws.Cells.Item(16, 5).Value = '  +23.32%  '

$ws.Cells.Item(17, 4).Value = '70.854.10'
$ws.Cells.Item(17, 5).Value = '  +6.09%  '

$ws.Cells.Item(18, 5).Value = '  +1.59%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.64'
$ws.Cells.Item(19, 5).Value = '  +9.64%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '525.16'
$ws.Cells.Item(20, 5).Value = '  +8.79%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.07'
$ws.Cells.Item(21, 5).Value = '  +3.64%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.53'
$ws.Cells.Item(22, 5).Value = '  +24.10%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.752'
$ws.Cells.Item(23, 5).Value = '  +10.12%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '88.32'
$ws.Cells.Item(24, 5).Value = '  +6.10%  '

$ws.Cells.Item(25, 5).Value = '  +12.18%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '13.62'
$ws.Cells.Item(26, 5).Value = '  +7.35%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.94'
$ws.Cells.Item(27, 5).Value = '  +6.10%  '

$ws.Cells.Item(28, 5).Value = '  +0.21%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.0000126'
$ws.Cells.Item(29, 5).Value = '  +36.90%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.52'
$ws.Cells.Item(30, 5).Value = '  +10.01%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.92'
$ws.Cells.Item(31, 5).Value = '  +13.49%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.98'
$ws.Cells.Item(32, 5).Value = '  +2.08%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '32.72'
$ws.Cells.Item(33, 5).Value = '  +17.37%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.116'
$ws.Cells.Item(34, 5).Value = '  +4.69%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.998'
$ws.Cells.Item(35, 5).Value = '  -0.20%  '

$ws.Cells.Item(36, 2).Value = 'Filecoin'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '6.22'
$ws.Cells.Item(36, 5).Value = '  +12.60%  '

$ws.Cells.Item(37, 2).Value = 'Mantle'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.06'
$ws.Cells.Item(37, 5).Value = '  +11.90%  '

$ws.Cells.Item(38, 5).Value = '  +10.70%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.20'
$ws.Cells.Item(39, 5).Value = '  +10.76%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.134'
$ws.Cells.Item(40, 5).Value = '  +8.94%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '51.69'
$ws.Cells.Item(41, 5).Value = '  +5.68%  '

$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '431.68'
$ws.Cells.Item(42, 5).Value = '  +17.05%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '3.163.24'
$ws.Cells.Item(43, 5).Value = '  +13.24%  '

$ws.Cells.Item(44, 2).Value = 'Cosmos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.91'
$ws.Cells.Item(44, 5).Value = '  +8.60%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '44.50'
$ws.Cells.Item(45, 5).Value = '  -4.78%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.79'
$ws.Cells.Item(46, 5).Value = '  +5.00%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0370'
$ws.Cells.Item(47, 5).Value = '  +8.58%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '27.77'
$ws.Cells.Item(48, 5).Value = '  +8.30%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '140.88'
$ws.Cells.Item(49, 5).Value = '  +4.52%  '

$ws.Cells.Item(50, 5).Value = '  +10.73%  '

